# Auto-generated Excel COM-interop script to apply scheduled-runner market-price updates
# to the Hyperion_Profits workbook (per-sheet leve profit calculations).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 142.85715
$ws.Range("I12").Value = 115
$ws.Range("K12").Value = 115
$ws.Range("M12").Value = 55
$ws.Range("H76").Value = 6070.0586
$ws.Range("I76").Value = 4213.1665
$ws.Range("J76").Value = 7082.909
$ws.Range("K76").Value = 4213.1665
$ws.Range("L76").Value = 7082.909
$ws.Range("M76").Value = -3898.1665
$ws.Range("N76").Value = -7712.909
$ws.Range("H79").Value = 6070.0586
$ws.Range("I79").Value = 4213.1665
$ws.Range("J79").Value = 7082.909
$ws.Range("K79").Value = 4213.1665
$ws.Range("L79").Value = 7082.909
$ws.Range("M79").Value = -3121.1665
$ws.Range("N79").Value = -9266.909
$ws.Range("H82").Value = 1454.909
$ws.Range("I82").Value = 1454.909
$ws.Range("K82").Value = 4364.727000000001
$ws.Range("M82").Value = -3958.727000000001
$ws.Range("H85").Value = 1454.909
$ws.Range("I85").Value = 1454.909
$ws.Range("K85").Value = 4364.727000000001
$ws.Range("M85").Value = -2960.727000000001
$ws.Range("H88").Value = 3671.7144
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3671.7144
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3671.7144
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4483.7144
$ws.Range("H91").Value = 3671.7144
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3671.7144
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3671.7144
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6479.7144
$ws.Range("H100").Value = 1263.8948
$ws.Range("I100").Value = 1151.8572
$ws.Range("J100").Value = 1577.6
$ws.Range("K100").Value = 1151.8572
$ws.Range("L100").Value = 1577.6
$ws.Range("M100").Value = -610.8571999999999
$ws.Range("N100").Value = -2659.6
$ws.Range("H112").Value = 8586.6875
$ws.Range("J112").Value = 9493
$ws.Range("L112").Value = 28479
$ws.Range("N112").Value = -30695
$ws.Range("H131").Value = 4337.773
$ws.Range("I131").Value = 2695.4
$ws.Range("J131").Value = 7857.143
$ws.Range("K131").Value = 8086.200000000001
$ws.Range("L131").Value = 23571.429
$ws.Range("M131").Value = -3046.200000000001
$ws.Range("N131").Value = -33651.429
$ws.Range("H132").Value = 27030510
$ws.Range("I132").Value = 33337162
$ws.Range("J132").Value = 1999.8572
$ws.Range("K132").Value = 100011486
$ws.Range("L132").Value = 5999.571599999999
$ws.Range("M132").Value = -100008956
$ws.Range("N132").Value = -11059.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4185.327
$ws.Range("I32").Value = 2299.513
$ws.Range("J32").Value = 9842.77
$ws.Range("K32").Value = 2299.513
$ws.Range("L32").Value = 9842.77
$ws.Range("M32").Value = -2012.513
$ws.Range("N32").Value = -10416.77
$ws.Range("H45").Value = 9620354
$ws.Range("I45").Value = 19233276
$ws.Range("K45").Value = 19233276
$ws.Range("M45").Value = -19232899
$ws.Range("H61").Value = 7411.3335
$ws.Range("I61").Value = 9053.6
$ws.Range("J61").Value = 3305.6667
$ws.Range("K61").Value = 9053.6
$ws.Range("L61").Value = 3305.6667
$ws.Range("M61").Value = -8841.6
$ws.Range("N61").Value = -3729.6667
$ws.Range("H74").Value = 38440.6
$ws.Range("I74").Value = 8707.6875
$ws.Range("K74").Value = 8707.6875
$ws.Range("M74").Value = -7833.6875
$ws.Range("H77").Value = 38440.6
$ws.Range("I77").Value = 8707.6875
$ws.Range("K77").Value = 43538.4375
$ws.Range("M77").Value = -39170.4375
$ws.Range("H122").Value = 775450.0600000001
$ws.Range("I122").Value = 3610.8667
$ws.Range("K122").Value = 10832.6001
$ws.Range("M122").Value = -8382.6001
$ws.Range("H132").Value = 16741.223
$ws.Range("I132").Value = 23157.4
$ws.Range("J132").Value = 8721
$ws.Range("K132").Value = 69472.20000000001
$ws.Range("L132").Value = 26163
$ws.Range("M132").Value = -66942.20000000001
$ws.Range("N132").Value = -31223
$ws.Range("H136").Value = 7411.3335
$ws.Range("I136").Value = 9053.6
$ws.Range("J136").Value = 3305.6667
$ws.Range("K136").Value = 27160.8
$ws.Range("L136").Value = 9917.000100000001
$ws.Range("M136").Value = -24610.8
$ws.Range("N136").Value = -15017.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1300
$ws.Range("I7").Value = 1700
$ws.Range("J7").Value = 700
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 700
$ws.Range("M7").Value = -1587
$ws.Range("N7").Value = -926
$ws.Range("H11").Value = 532.5
$ws.Range("J11").Value = 740
$ws.Range("L11").Value = 740
$ws.Range("N11").Value = -1020
$ws.Range("H76").Value = 10143
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 10143
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value = 16550
$ws.Range("I82").Value = 17066.666
$ws.Range("J82").Value = 15000
$ws.Range("K82").Value = 17066.666
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -16683.666
$ws.Range("N82").Value = -15766
$ws.Range("H85").Value = 16550
$ws.Range("I85").Value = 17066.666
$ws.Range("J85").Value = 15000
$ws.Range("K85").Value = 17066.666
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -15740.666
$ws.Range("N85").Value = -17652
$ws.Range("H134").Value = 11391.531
$ws.Range("I134").Value = 11946.218
$ws.Range("K134").Value = 35838.654
$ws.Range("M134").Value = -33303.654

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31443.555
$ws.Range("I31").Value = 11007.546
$ws.Range("K31").Value = 11007.546
$ws.Range("M31").Value = -10712.546
$ws.Range("H34").Value = 31443.555
$ws.Range("I34").Value = 11007.546
$ws.Range("K34").Value = 11007.546
$ws.Range("M34").Value = -10805.546
$ws.Range("H36").Value = 2048
$ws.Range("I36").Value = 2048
$ws.Range("K36").Value = 2048
$ws.Range("M36").Value = -1660
$ws.Range("H40").Value = 2048
$ws.Range("I40").Value = 2048
$ws.Range("K40").Value = 2048
$ws.Range("M40").Value = -1888
$ws.Range("H105").Value = 2191.6428
$ws.Range("I105").Value = 2098.5833
$ws.Range("K105").Value = 2098.5833
$ws.Range("M105").Value = -351.5832999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1300
$ws.Range("I5").Value = 1350
$ws.Range("K5").Value = 4050
$ws.Range("M5").Value = -3938
$ws.Range("H33").Value = 6308.125
$ws.Range("I33").Value = 2695.7368
$ws.Range("J33").Value = 20035.2
$ws.Range("K33").Value = 16174.4208
$ws.Range("L33").Value = 120211.2
$ws.Range("M33").Value = -15891.4208
$ws.Range("N33").Value = -120777.2
$ws.Range("H37").Value = 68316
$ws.Range("J37").Value = 68316
$ws.Range("L37").Value = 204948
$ws.Range("N37").Value = -205172
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H135").Value = 1300
$ws.Range("I135").Value = 1350
$ws.Range("K135").Value = 12150
$ws.Range("M135").Value = -9615
$ws.Range("H140").Value = 2359.0588
$ws.Range("I140").Value = 2194
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 6582
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -1402
$ws.Range("N140").Value = -25360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 890645.2
$ws.Range("I122").Value = 989516.9
$ws.Range("K122").Value = 2968550.7
$ws.Range("M122").Value = -2966100.7
$ws.Range("H132").Value = 10194.482
$ws.Range("I132").Value = 7090.9565
$ws.Range("J132").Value = 22091.334
$ws.Range("K132").Value = 21272.8695
$ws.Range("L132").Value = 66274.00199999999
$ws.Range("M132").Value = -18742.8695
$ws.Range("N132").Value = -71334.00199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 25643822
$ws.Range("I93").Value = 41669270
$ws.Range("J93").Value = 3099.2
$ws.Range("K93").Value = 41669270
$ws.Range("L93").Value = 3099.2
$ws.Range("M93").Value = -41668022
$ws.Range("N93").Value = -5595.2
$ws.Range("H136").Value = 40847.035
$ws.Range("J136").Value = 7432.25
$ws.Range("L136").Value = 22296.75
$ws.Range("N136").Value = -27396.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3442.875
$ws.Range("I122").Value = 3637.125
$ws.Range("J122").Value = 3248.625
$ws.Range("K122").Value = 10911.375
$ws.Range("L122").Value = 9745.875
$ws.Range("M122").Value = -8461.375
$ws.Range("N122").Value = -14645.875
$ws.Range("H132").Value = 84258940
$ws.Range("I132").Value = 166684530
$ws.Range("K132").Value = 500053590
$ws.Range("M132").Value = -500051060
